# Added mtcars data set: the summary table grows from 4 cylinder-count
# groups (rows 5-8) to 5 groups (rows 5-9), and the empty footer row
# (formerly row 9) shifts down to row 10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push the existing "totals" row (row 9, style 10) down to row 10.
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

# 2) Turn row 9 into a data row matching the formatting of rows 5-8.
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)

$ws.Range("A9:E9").ClearContents()

# 3) Write the new mtcars-derived values.
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 91
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 2.14
$ws.Range("E5").ClearContents()

$ws.Range("A6").Value = 10
$ws.Range("B6").Value = 81.8
$ws.Range("C6").Value = 21.87235698318771
$ws.Range("D6").Value = 2.3003
$ws.Range("E6").Value = 0.5982073312080948

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 131.6666666666667
$ws.Range("C7").Value = 37.52776749732568
$ws.Range("D7").Value = 2.755
$ws.Range("E7").Value = 0.1281600561797629

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = 115.25
$ws.Range("C8").Value = 9.178779875342908
$ws.Range("D8").Value = 3.38875
$ws.Range("E8").Value = 0.1162163929916946

$ws.Range("A9").Value = 14
$ws.Range("B9").Value = 209.2142857142857
$ws.Range("C9").Value = 50.97688551827051
$ws.Range("D9").Value = 3.999214285714287
$ws.Range("E9").Value = 0.7594047444769265
